$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New test case rows (Booking / Discount code)
# ---------------------------------------------------------------------------

# Row 2 - Valid Discount code
$ws.Rows.Item(2).RowHeight = 91
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Booking"
$ws.Range("C2").Value = "Valid Discount code " + [char]0x2013 + " user can proceed with payment"
$ws.Range("D2").Value = "User is on passenger details page."
$ws.Range("E2").Value = "1. Enter a valid email, first and last name.`n2. Fill other required fields correctly.`n3. Get a correct discount code and set it to the Discount code field`n4. Click continue."
$ws.Range("F2").Value = "System allows user to proceed with payment and applies all data including discount code"
$ws.Range("G2").Value = "High"
$ws.Range("H2").Value = "Positive / Validation"
$ws.Range("J2").Value = "To_Be_Automated"

# Row 3 - Discount code field is empty
$ws.Rows.Item(3).RowHeight = 76.1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Booking"
$ws.Range("C3").Value = "Discount code field is empty"
$ws.Range("D3").Value = "User is on passenger details page."
$ws.Range("E3").Value = "1. Enter a valid email, first and last name.`n2. Fill other required fields correctly.`n3. Leave Discount code field empty`n4. Click continue."
$ws.Range("F3").Value = "System displays a validation error indicating the discount code is invalid and does not proceed."
$ws.Range("G3").Value = "High"
$ws.Range("H3").Value = "Negative / Validation"
$ws.Range("J3").Value = "To_Be_Automated"

# Row 4 - Invalid discount code value
$ws.Rows.Item(4).RowHeight = 76.1
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Booking"
$ws.Range("C4").Value = "Invalid discount code value"
$ws.Range("D4").Value = "User is on passenger details page."
$ws.Range("E4").Value = "1. Enter a valid email, first and last name.`n2. Fill other required fields correctly.`n3. Set an invalid discount code`n4. Click continue."
$ws.Range("F4").Value = "System displays a validation error indicating the discount code is invalid and does not proceed."
$ws.Range("G4").Value = "High"
$ws.Range("H4").Value = "Negative / Validation"
$ws.Range("J4").Value = "To_Be_Automated"

# ---------------------------------------------------------------------------
# E3/E4 and J3/J4 use a dedicated wrapped style (new font + new cell format)
# identical to the one used for the multi-line Test Steps / Notes cells.
# Build it once on the "3" row, then clone the resolved style onto row 4 via
# a format-only paste so we don't create duplicate / orphaned cell formats.
# ---------------------------------------------------------------------------
$ws.Range("E3").WrapText = $true
$ws.Range("E3").Font.ThemeColor = 1
$ws.Range("E3").Font.Name = "Calibri"
$ws.Range("E3").Font.Size = 11
$ws.Range("E3").Font.Family = 2
$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("J3").WrapText = $true
$ws.Range("J3").HorizontalAlignment = -4108
$ws.Range("J3").VerticalAlignment = -4108
$ws.Range("J3").Font.ThemeColor = 1
$ws.Range("J3").Font.Name = "Calibri"
$ws.Range("J3").Font.Size = 11
$ws.Range("J3").Font.Family = 2
$ws.Range("J3").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Update the active selection to E4, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()
